# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) counts for rows 2-25 (column G), replacing the old Strike# counts
$kValues = @(3, 6, 5, 6, 1, 3, 5, 5, 5, 7, 7, 7, 3, 7, 9, 5, 10, 6, 9, 6, 6, 5, 2, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
